# Generate Report for Handback
# This script reproduces the "handback report" update for the
# 2a4db539-2e4d-40af-b43d-b1943adfb440.md row (row 6) on both the
# zh-cn and de-de localization-status sheets, plus the related column
# width widening for the "Error Detail" column.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb654d065f38485e679122520789daf8c6f35562/e2e/2a4db539-2e4d-40af-b43d-b1943adfb440.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ea258442fa1151debbb1203d22d97c3b870234b/e2e/2a4db539-2e4d-40af-b43d-b1943adfb440.md."
$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ea258442fa1151debbb1203d22d97c3b870234b/e2e/2a4db539-2e4d-40af-b43d-b1943adfb440.md"
$mdDisplay = "2a4db539-2e4d-40af-b43d-b1943adfb440.md"

function Update-HandbackSheet($SheetName, $XlfName, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Widen the "Error Detail" column (P, the 16th column) to fit the
    # new long error message (COM ColumnWidth 39.17 -> OOXML width 40).
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # Preserve the existing hyperlinks that live after row 6 / column A
    # so we can re-create them after inserting the new one -- this keeps
    # the hyperlink / relationship ordering identical to a natural
    # "insert in the middle" edit.
    $laterRefs = @()
    $laterTargets = @()
    $laterDisplays = @()
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$A$7' -or $addr -eq '$A$8') {
            $laterRefs += ($addr -replace '\$', '')
            $laterTargets += $hl.Address()
            $laterDisplays += $hl.TextToDisplay()
        }
    }

    for ($i = 0; $i -lt $laterRefs.Count; $i++) {
        $ref = $laterRefs[$i]
        foreach ($hl in $ws.Hyperlinks) {
            $addr = $hl.Range.Address() -replace '\$', ''
            if ($addr -eq $ref) {
                $hl.Delete()
            }
        }
    }

    # Fill in the "Latest Target File", "Latest Handback File",
    # "Latest Handback DateTime" and "Error Detail" columns for row 6.
    $ws.Range("J6").Value = $XlfName
    $ws.Range("K6").Value = $HandbackDateTime
    $ws.Range("P6").Value = $errorMessage

    $ws.Range("I6").Value = $mdDisplay
    $ws.Hyperlinks.Add($ws.Range("I6"), $latestMdUrl, "", "", $mdDisplay) | Out-Null
    $ws.Range("I6").Style = "HyperLink"

    # Re-create the hyperlinks for A7 / A8 (now shifted after the new
    # relationship that was inserted for I6).
    for ($i = 0; $i -lt $laterRefs.Count; $i++) {
        $ws.Hyperlinks.Add($ws.Range($laterRefs[$i]), $laterTargets[$i], "", "", $laterDisplays[$i]) | Out-Null
    }
}

Update-HandbackSheet "zh-cn" "2a4db539-2e4d-40af-b43d-b1943adfb440.4bb7b0a80ecaa41d71f1f06a3f7c698760343090.zh-cn.xlf" "2016-09-01 12:48:23"
Update-HandbackSheet "de-de" "2a4db539-2e4d-40af-b43d-b1943adfb440.4bb7b0a80ecaa41d71f1f06a3f7c698760343090.de-de.xlf" "2016-09-01 12:48:30"

$wb.Save()
